# Applies the "Added the creation of the user and the generation of the token"
# commit to the workbook: fixes the Data-sheet template so the Yes/No
# correct-answer flags live in E:I (matching the "answer_N_value" headers)
# and the sample answer text lives in J:N (matching the "answer_N" headers),
# updates the helper formulas and the Yes/No data-validation range to match,
# switches the sample question type to "Single Selection", and leaves the
# selection where the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# --- B2: sample question type Multi Selection -> Single Selection ---
# (B2:B10 share the same quote-prefixed "text" style; re-apply it with a
# formats-only paste from a neighbouring cell so the style index survives
# the value change.)
$ws.Range("B2").Value = "Single Selection"
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Swap the answer-text columns (E:I) with the Yes/No columns (J:N) ---
# The template previously had the correct-answer Yes/No flags in J:N and the
# placeholder answer text in E:I; the edit moves the Yes/No flags into E:I
# (matching the E:I "answer_N_value" headers) and the answer text into J:N.
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = "Yes"
$ws.Range("H2").Value = "No"
$ws.Range("I2").Value = "No"
$ws.Range("J2").Value = "answer_1_value1"
$ws.Range("K2").Value = "answer_2_value1"
$ws.Range("L2").Value = "answer_3_value1"
$ws.Range("M2").Value = "answer_4_value1"
$ws.Range("N2").Value = "answer_5_value1"

$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = "No"
$ws.Range("H3").Value = "Yes"
$ws.Range("I3").Value = "No"
$ws.Range("J3").Value = "answer_1_value1"
$ws.Range("K3").Value = "answer_2_value1"
$ws.Range("L3").Value = "answer_3_value1"
$ws.Range("M3").Value = "answer_4_value1"
$ws.Range("N3").Value = "answer_5_value1"

$ws.Range("E4").Value = "Yes"
$ws.Range("F4").Value = "No"
$ws.Range("G4").Value = "Yes"
$ws.Range("H4").Value = "Yes"
$ws.Range("I4").Value = "No"
$ws.Range("J4").Value = "answer_1_value1"
$ws.Range("K4").Value = "answer_2_value1"
$ws.Range("L4").Value = "answer_3_value1"
$ws.Range("M4").Value = "answer_4_value1"
$ws.Range("N4").Value = "answer_5_value1"

$ws.Range("E5").Value = "Yes"
$ws.Range("F5").Value = "Yes"
$ws.Range("G5").Value = "Yes"
$ws.Range("H5").Value = "No"
$ws.Range("I5").Value = "No"
$ws.Range("J5").Value = "answer_1_value1"
$ws.Range("K5").Value = "answer_2_value1"
$ws.Range("L5").Value = "answer_3_value1"
$ws.Range("M5").Value = "answer_4_value1"
$ws.Range("N5").Value = "answer_5_value1"

$ws.Range("E6").Value = "Yes"
$ws.Range("F6").Value = "No"
$ws.Range("G6").Value = "Yes"
$ws.Range("H6").Value = "No"
$ws.Range("I6").Value = "Yes"
$ws.Range("J6").Value = "answer_1_value1"
$ws.Range("K6").Value = "answer_2_value1"
$ws.Range("L6").Value = "answer_3_value1"
$ws.Range("M6").Value = "answer_4_value1"
$ws.Range("N6").Value = "answer_5_value1"

$ws.Range("E7").Value = "Yes"
$ws.Range("F7").Value = "Yes"
$ws.Range("G7").Value = "Yes"
$ws.Range("H7").Value = "No"
$ws.Range("I7").Value = "No"
$ws.Range("J7").Value = "answer_1_value1"
$ws.Range("K7").Value = "answer_2_value1"
$ws.Range("L7").Value = "answer_3_value1"
$ws.Range("M7").Value = "answer_4_value1"
$ws.Range("N7").Value = "answer_5_value1"

$ws.Range("E8").Value = "Yes"
$ws.Range("F8").Value = "Yes"
$ws.Range("G8").Value = "Yes"
$ws.Range("H8").Value = "No"
$ws.Range("I8").Value = "No"
$ws.Range("J8").Value = "answer_1_value1"
$ws.Range("K8").Value = "answer_2_value1"
$ws.Range("L8").Value = "answer_3_value1"
$ws.Range("M8").Value = "answer_4_value1"
$ws.Range("N8").Value = "answer_5_value1"

$ws.Range("E9").Value = "Yes"
$ws.Range("F9").Value = "Yes"
$ws.Range("G9").Value = "Yes"
$ws.Range("H9").Value = "No"
$ws.Range("I9").Value = "No"
$ws.Range("J9").Value = "answer_1_value1"
$ws.Range("K9").Value = "answer_2_value1"
$ws.Range("L9").Value = "answer_3_value1"
$ws.Range("M9").Value = "answer_4_value1"
$ws.Range("N9").Value = "answer_5_value1"

$ws.Range("E10").Value = "Yes"
$ws.Range("F10").Value = "Yes"
$ws.Range("G10").Value = "Yes"
$ws.Range("H10").Value = "No"
$ws.Range("I10").Value = "No"
$ws.Range("J10").Value = "answer_1_value1"
$ws.Range("K10").Value = "answer_2_value1"
$ws.Range("L10").Value = "answer_3_value1"
$ws.Range("M10").Value = "answer_4_value1"
$ws.Range("N10").Value = "answer_5_value1"

# --- Repoint the helper formulas (P:U) at the swapped columns ---
$ws.Range("P2").Formula = "=IF(E2=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q2").Formula = "=IF(F2=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R2").Formula = "=IF(G2=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S2").Formula = "=IF(H2=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T2").Formula = "=IF(I2=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U2").Formula = "=IF(N2=""Yes"",I$1,"""")"

$ws.Range("P3").Formula = "=IF(E3=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q3").Formula = "=IF(F3=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R3").Formula = "=IF(G3=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S3").Formula = "=IF(H3=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T3").Formula = "=IF(I3=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U3").Formula = "=IF(N3=""Yes"",I$1,"""")"

$ws.Range("P4").Formula = "=IF(E4=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q4").Formula = "=IF(F4=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R4").Formula = "=IF(G4=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S4").Formula = "=IF(H4=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T4").Formula = "=IF(I4=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U4").Formula = "=IF(N4=""Yes"",I$1,"""")"

$ws.Range("P5").Formula = "=IF(E5=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q5").Formula = "=IF(F5=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R5").Formula = "=IF(G5=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S5").Formula = "=IF(H5=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T5").Formula = "=IF(I5=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U5").Formula = "=IF(N5=""Yes"",I$1,"""")"

$ws.Range("P6").Formula = "=IF(E6=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q6").Formula = "=IF(F6=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R6").Formula = "=IF(G6=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S6").Formula = "=IF(H6=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T6").Formula = "=IF(I6=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U6").Formula = "=IF(N6=""Yes"",I$1,"""")"

$ws.Range("P7").Formula = "=IF(E7=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q7").Formula = "=IF(F7=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R7").Formula = "=IF(G7=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S7").Formula = "=IF(H7=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T7").Formula = "=IF(I7=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U7").Formula = "=IF(N7=""Yes"",I$1,"""")"

$ws.Range("P8").Formula = "=IF(E8=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q8").Formula = "=IF(F8=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R8").Formula = "=IF(G8=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S8").Formula = "=IF(H8=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T8").Formula = "=IF(I8=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U8").Formula = "=IF(N8=""Yes"",I$1,"""")"

$ws.Range("P9").Formula = "=IF(E9=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q9").Formula = "=IF(F9=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R9").Formula = "=IF(G9=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S9").Formula = "=IF(H9=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T9").Formula = "=IF(I9=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U9").Formula = "=IF(N9=""Yes"",I$1,"""")"

$ws.Range("P10").Formula = "=IF(E10=""Yes"",CONCATENATE(J$1,"",""),"""")"
$ws.Range("Q10").Formula = "=IF(F10=""Yes"",CONCATENATE(K$1,"",""),"""")"
$ws.Range("R10").Formula = "=IF(G10=""Yes"",CONCATENATE(L$1,"",""),"""")"
$ws.Range("S10").Formula = "=IF(H10=""Yes"",CONCATENATE(M$1,"",""),"""")"
$ws.Range("T10").Formula = "=IF(I10=""Yes"",CONCATENATE(N$1,"",""),"""")"
$ws.Range("U10").Formula = "=IF(N10=""Yes"",I$1,"""")"


# --- Move the Yes/No list validation from J2:N10 to E2:I10 ---
$ws.Range("E2:I10").Validation.Add(3, 1, 1, """Yes,No""")
$ws.Range("J2:N10").Validation.Delete()

# --- Restore the cursor position left behind when the author saved ---
$ws.Range("F19").Select()
